# Weekly update: insert a new Naranja price record as the latest week's
# reading, pushing all prior weekly records down by one row
# (Fruta / hortaliza, semanal).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 46 - this shifts existing rows 46:92 down to 47:93
# and updates the used-range dimension automatically.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with this week's record. Most of the
# columns mirror the row immediately below (same market / product taxonomy),
# only the date, variety, volume and origin region are new for this entry.
$ws.Range("A46").Value = 1
$ws.Range("B46").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C46").Value = "Arica y Parinacota"
$ws.Range("D46").Value = 44763
$ws.Range("E46").Value = 15
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100102
$ws.Range("H46").Value = "Cítricos"
$ws.Range("I46").Value = 100102005
$ws.Range("J46").Value = "Naranja"
$ws.Range("K46").Value = "Thompson"
$ws.Range("L46").Value = "Segunda"
$ws.Range("M46").Value = 270
$ws.Range("N46").Value = 700
$ws.Range("O46").Value = 750
$ws.Range("P46").Value = 725
$ws.Range("Q46").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R46").Value = "Región de O'Higgins"
$ws.Range("S46").Value = 725
$ws.Range("T46").Value = 1
